$wb = $excel.ActiveWorkbook

# ---- Richness (sheet1) ----
$ws = $wb.Worksheets.Item("Richness")

# Comparison labels (shared text changed across all sheets)
$ws.Range("A2").Value = "healthy vs Grouppre_ltx"
$ws.Range("A3").Value = "healthy vs pre_ltx - CZ vs NO"
$ws.Range("A4").Value = "healthy vs Grouppre_ltx:CountryNO"
$ws.Range("A8").Value = "healthy vs Grouppost_ltx"
$ws.Range("A9").Value = "healthy vs post_ltx - CZ vs NO"
$ws.Range("A10").Value = "healthy vs Grouppost_ltx:CountryNO"

# Numeric / sig recalculated values
$ws.Range("B2").Value = -56.2297874111623
$ws.Range("C2").Value = 11.7376584507647
$ws.Range("D2").Value = 195.570180915709
$ws.Range("E2").Value = -4.79054554594736
$ws.Range("F2").Value = 0.00000328100031952592
$ws.Range("G2").Value = 0.0000295290028757333
$ws.Range("B3").Value = 14.1930992929996
$ws.Range("C3").Value = 10.3411147337777
$ws.Range("D3").Value = 194.537858584009
$ws.Range("E3").Value = 1.37249219821921
$ws.Range("F3").Value = 0.171490557299805
$ws.Range("G3").Value = 0.255969836828398
$ws.Range("B4").Value = -7.48936434332061
$ws.Range("C4").Value = 15.1691326723597
$ws.Range("D4").Value = 194.358194950027
$ws.Range("E4").Value = -0.493723965969872
$ws.Range("F4").Value = 0.622058607059278
$ws.Range("G4").Value = 0.622058607059278
$ws.Range("G6").Value = 0.622058607059278
$ws.Range("G7").Value = 0.202957793373846
$ws.Range("B8").Value = -27.299148916624
$ws.Range("C8").Value = 8.69854137378728
$ws.Range("D8").Value = 250.458335818147
$ws.Range("E8").Value = -3.13835938044611
$ws.Range("F8").Value = 0.0019024383230054
$ws.Range("G8").Value = 0.00856097245352429
$ws.Range("B9").Value = 14.3603645603166
$ws.Range("C9").Value = 11.1533109482603
$ws.Range("D9").Value = 251.728796047421
$ws.Range("E9").Value = 1.28754274196547
$ws.Range("F9").Value = 0.199087650866532
$ws.Range("G9").Value = 0.255969836828398
$ws.Range("B10").Value = -31.6814031950151
$ws.Range("C10").Value = 14.19496184593
$ws.Range("D10").Value = 247.530840100213
$ws.Range("E10").Value = -2.23187660092929
$ws.Range("F10").Value = 0.0265186605242804
$ws.Range("G10").Value = 0.0596669861796309

# ---- Shannon (sheet2) ----
$ws = $wb.Worksheets.Item("Shannon")

# Comparison labels (shared text changed across all sheets)
$ws.Range("A2").Value = "healthy vs Grouppre_ltx"
$ws.Range("A3").Value = "healthy vs pre_ltx - CZ vs NO"
$ws.Range("A4").Value = "healthy vs Grouppre_ltx:CountryNO"
$ws.Range("A8").Value = "healthy vs Grouppost_ltx"
$ws.Range("A9").Value = "healthy vs post_ltx - CZ vs NO"
$ws.Range("A10").Value = "healthy vs Grouppost_ltx:CountryNO"

# Numeric / sig recalculated values
$ws.Range("B2").Value = -0.584638716137957
$ws.Range("C2").Value = 0.143155036381466
$ws.Range("D2").Value = 197.587543805541
$ws.Range("E2").Value = -4.08395492688128
$ws.Range("F2").Value = 0.0000642923957847682
$ws.Range("G2").Value = 0.000578631562062914
$ws.Range("B3").Value = 0.0146781534706523
$ws.Range("C3").Value = 0.126064800286197
$ws.Range("D3").Value = 196.213354178243
$ws.Range("E3").Value = 0.116433401213736
$ws.Range("F3").Value = 0.907428148748997
$ws.Range("G3").Value = 0.909474207145245
$ws.Range("B4").Value = 0.0807676638695282
$ws.Range("C4").Value = 0.184906372790758
$ws.Range("D4").Value = 195.965218122722
$ws.Range("E4").Value = 0.436803029828105
$ws.Range("F4").Value = 0.662735115053663
$ws.Range("G4").Value = 0.852088005068996
$ws.Range("G5").Value = 0.0564792646200416
$ws.Range("H5").Value = ""
$ws.Range("G6").Value = 0.838557559871001
$ws.Range("G7").Value = 0.0564792646200416
$ws.Range("H7").Value = ""
$ws.Range("B8").Value = -0.219739423568967
$ws.Range("C8").Value = 0.102071461452386
$ws.Range("D8").Value = 247.585634183915
$ws.Range("E8").Value = -2.1527998173267
$ws.Range("F8").Value = 0.0323009116336897
$ws.Range("G8").Value = 0.0581416409406414
$ws.Range("H8").Value = ""
$ws.Range("B9").Value = 0.0148944623604509
$ws.Range("C9").Value = 0.130862581417962
$ws.Range("D9").Value = 248.75891335263
$ws.Range("E9").Value = 0.113817580236168
$ws.Range("F9").Value = 0.909474207145245
$ws.Range("G9").Value = 0.909474207145245
$ws.Range("H9").Value = ""
$ws.Range("B10").Value = -0.375725890122552
$ws.Range("C10").Value = 0.166609248148072
$ws.Range("D10").Value = 244.90353833122
$ws.Range("E10").Value = -2.25513225885654
$ws.Range("F10").Value = 0.0250086183896467
$ws.Range("G10").Value = 0.0564792646200416
$ws.Range("H10").Value = ""

# ---- Simpson (sheet3) ----
$ws = $wb.Worksheets.Item("Simpson")

# Comparison labels (shared text changed across all sheets)
$ws.Range("A2").Value = "healthy vs Grouppre_ltx"
$ws.Range("A3").Value = "healthy vs pre_ltx - CZ vs NO"
$ws.Range("A4").Value = "healthy vs Grouppre_ltx:CountryNO"
$ws.Range("A8").Value = "healthy vs Grouppost_ltx"
$ws.Range("A9").Value = "healthy vs post_ltx - CZ vs NO"
$ws.Range("A10").Value = "healthy vs Grouppost_ltx:CountryNO"

# Numeric / sig recalculated values
$ws.Range("B2").Value = -0.058693137605951
$ws.Range("C2").Value = 0.0191445452131798
$ws.Range("D2").Value = 199.8531786722
$ws.Range("E2").Value = -3.06578907737878
$ws.Range("F2").Value = 0.00247109029994167
$ws.Range("G2").Value = 0.0222398126994751
$ws.Range("B3").Value = -0.00725094013991043
$ws.Range("C3").Value = 0.0168510078148344
$ws.Range("D3").Value = 198.127375767928
$ws.Range("E3").Value = -0.430297120480072
$ws.Range("F3").Value = 0.667447205149642
$ws.Range("G3").Value = 0.737867463605542
$ws.Range("B4").Value = 0.0200982842190763
$ws.Range("C4").Value = 0.0247141535954386
$ws.Range("D4").Value = 197.803443913675
$ws.Range("E4").Value = 0.813229720429746
$ws.Range("F4").Value = 0.417064141580538
$ws.Range("G4").Value = 0.625596212370807
$ws.Range("G5").Value = 0.45471288265286
$ws.Range("G6").Value = 0.737867463605542
$ws.Range("G7").Value = 0.362553208487713
$ws.Range("B8").Value = -0.0270796696678986
$ws.Range("C8").Value = 0.0169835936225165
$ws.Range("D8").Value = 246.082173765515
$ws.Range("E8").Value = -1.59446052877743
$ws.Range("F8").Value = 0.112116439286781
$ws.Range("G8").Value = 0.362553208487713
$ws.Range("B9").Value = -0.00729458643870464
$ws.Range("C9").Value = 0.0217712802387507
$ws.Range("D9").Value = 247.135382698276
$ws.Range("E9").Value = -0.335055465673581
$ws.Range("F9").Value = 0.737867463605542
$ws.Range("G9").Value = 0.737867463605542
$ws.Range("B10").Value = -0.0317991358794888
$ws.Range("C10").Value = 0.02773029605283
$ws.Range("D10").Value = 243.697487853012
$ws.Range("E10").Value = 1.14672904389146
$ws.Range("F10").Value = 0.252618268140478
$ws.Range("G10").Value = 0.45471288265286

# ---- Pielou (sheet4) ----
$ws = $wb.Worksheets.Item("Pielou")

# Comparison labels (shared text changed across all sheets)
$ws.Range("A2").Value = "healthy vs Grouppre_ltx"
$ws.Range("A3").Value = "healthy vs pre_ltx - CZ vs NO"
$ws.Range("A4").Value = "healthy vs Grouppre_ltx:CountryNO"
$ws.Range("A8").Value = "healthy vs Grouppost_ltx"
$ws.Range("A9").Value = "healthy vs post_ltx - CZ vs NO"
$ws.Range("A10").Value = "healthy vs Grouppost_ltx:CountryNO"

# Numeric / sig recalculated values
$ws.Range("B2").Value = -0.0569586378083008
$ws.Range("C2").Value = 0.0201379137971631
$ws.Range("D2").Value = 201.129821433303
$ws.Range("E2").Value = -2.82842792863304
$ws.Range("F2").Value = 0.00515117600663822
$ws.Range("G2").Value = 0.046360584059744
$ws.Range("B3").Value = -0.00790346060303697
$ws.Range("C3").Value = 0.0177159654026436
$ws.Range("D3").Value = 199.016528896297
$ws.Range("E3").Value = -0.44612079688627
$ws.Range("F3").Value = 0.655994977824663
$ws.Range("G3").Value = 0.78097352770903
$ws.Range("B4").Value = 0.014685728771077
$ws.Range("C4").Value = 0.0259800634037992
$ws.Range("D4").Value = 198.602151099978
$ws.Range("E4").Value = 0.565269165930113
$ws.Range("F4").Value = 0.572528702844586
$ws.Range("G4").Value = 0.78097352770903
$ws.Range("G5").Value = 0.234319652628036
$ws.Range("G7").Value = 0.202339306421313
$ws.Range("B8").Value = -0.0178562322308237
$ws.Range("C8").Value = 0.0158041262922856
$ws.Range("D8").Value = 247.076162786932
$ws.Range("E8").Value = -1.12984621234897
$ws.Range("F8").Value = 0.259636697000521
$ws.Range("G8").Value = 0.467346054600939
$ws.Range("B9").Value = -0.0079773077315203
$ws.Range("C9").Value = 0.0202664850400594
$ws.Range("D9").Value = 248.451075779612
$ws.Range("E9").Value = -0.393620685370556
$ws.Range("F9").Value = 0.694198691296915
$ws.Range("G9").Value = 0.78097352770903
$ws.Range("H9").Value = ""
$ws.Range("B10").Value = -0.0471222956091828
$ws.Range("C10").Value = 0.0257834023145581
$ws.Range("D10").Value = 243.881388465963
$ws.Range("E10").Value = -1.82762131367651
$ws.Range("F10").Value = 0.0688282457285189
$ws.Range("G10").Value = 0.206484737185557
